# update database and change read_price algorithm
# Inserts 5 new historical quarterly columns (D:H) before the existing
# data block (which shifts to I:M), fills them with the new data pulled
# from the updated source, and touches up a few cosmetic knobs (theme
# accent1/accent5 swap, new column widths) that came along with the
# refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Make room: insert 5 fresh columns at D, pushing the old D:H block
#    (existing 1400/09 .. 1401/09 data) out to I:M.
# ---------------------------------------------------------------------
$ws.Range("D1:H1").EntireColumn.Insert()

# ---------------------------------------------------------------------
# 2. Column widths. D:E and G:I keep the "28" width, F keeps "29" (same
#    repeating pattern as the original D:H block).
# ---------------------------------------------------------------------
$ws.Columns("D:E").ColumnWidth = 27.1667
$ws.Columns("F").ColumnWidth = 28.1667
$ws.Columns("G:H").ColumnWidth = 27.1667

# ---------------------------------------------------------------------
# 3. Period headers (row 8) and publish-date headers (row 9) for the
#    5 new columns.
# ---------------------------------------------------------------------
$ws.Cells.Item(8, 4).Value = "6 ماهه منتهی به 1399/06"
$ws.Cells.Item(8, 5).Value = "9 ماهه منتهی به 1399/09"
$ws.Cells.Item(8, 6).Value = "12 ماهه منتهی به 1399/12"
$ws.Cells.Item(8, 7).Value = "3 ماهه منتهی به 1400/03"
$ws.Cells.Item(8, 8).Value = "6 ماهه منتهی به 1400/06"

$ws.Cells.Item(9, 4).Value = "1400-09-30 (4)"
$ws.Cells.Item(9, 5).Value = "1400-10-30 (2)"
$ws.Cells.Item(9, 6).Value = "1401-04-15 (8)"
$ws.Cells.Item(9, 7).Value = "1401-04-29 (2)"
$ws.Cells.Item(9, 8).Value = "1401-09-14 (4)"

# ---------------------------------------------------------------------
# 4. Financial data for the 5 new columns, row by row.
# ---------------------------------------------------------------------
$data = @{
    11 = @(223480, 377842, 589013, 214103, 310000)
    12 = @(-141547, -244387, -419579, -172248, -245353)
    13 = @(81933, 133455, 169434, 41855, 64647)
    14 = @(-3363, -5511, -9208, -3204, -5134)
    15 = @("-", "-", "-", "-", "-")
    16 = @(-3824, 2842, 2831, -1820, -4248)
    17 = @(74746, 130787, 163056, 36831, 55266)
    18 = @(-7450, -12521, -15192, -4382, -6584)
    19 = @(531, 438, 20620, -1500, -1618)
    20 = @(67828, 118704, 168484, 30949, 47064)
    21 = @(-5536, -4875, -16214, "-", -2672)
    22 = @(62291, 113829, 152270, 30949, 44391)
    23 = @("-", "-", "-", "-", "-")
    24 = @(62291, 113829, 152270, 30949, 44391)
    25 = @(0, 0, 0, 0, 0)
    26 = @(77442, 68194, 66342, 64295, 129133)
    27 = @(0, 0, 0, 0, 0)
}

foreach ($row in 11..27) {
    $vals = $data[$row]
    if ($null -eq $vals) { continue }
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($row, 4 + $i).Value = $vals[$i]
    }
}

# ---------------------------------------------------------------------
# 5. Theme refresh: accent1 and accent5 swapped places.
# ---------------------------------------------------------------------
$tcs = $wb.Theme.ThemeColorScheme
$tcs.Colors(5).RGB = 12874308   # accent1 <- 4472C4
$tcs.Colors(9).RGB = 13998939   # accent5 <- 5B9BD5

Write-Output "edit applied"
